$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Structural edits first (row insert / reshuffle), values and
# formatting are filled in afterwards.
# ============================================================

# Row 27: hours 3 -> 3.5, append a note to the log entry, grow the row height
$ws.Range("B27").Value = 3.5
$existing27 = $ws.Range("D27").Formula
$ws.Range("D27").Value = $existing27 + "`nWeek 5: started videos"
$ws.Rows.Item(27).RowHeight = 75

# New row 28: next day's date (Insert copies formatting down from row 27, so
# A28 picks up the date style; clear the stray D28 cell it also creates).
# This insert also pushes the old rows 33/34/36/38 down to 34/35/37/39.
$ws.Rows.Item(28).Insert()
$ws.Range("D28").Clear()
$ws.Range("A28").Value = 43521

# The note that used to be row 38 is now (after the shift above) sitting at
# row 39; push it further down to row 42 by inserting three more blank rows
# right before it. Rows 34/35/37 (old 33/34/36) are unaffected since they're
# above the insertion point.
$ws.Range("A39:A41").EntireRow.Insert()

# ============================================================
# Fill in the brand-new text cells, in the same order the source
# workbook's shared-string table lists them so the indices line up.
# ============================================================
$ws.Range("D27").Copy()

$ws.Range("D32").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D32").Value = "Issues/Loose Ends:"

$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value = "// Instead of comparing all values`n        // it may make sense to use .equals()`n        // TODO review .equals recommendations http://docs.jboss.org/hibernate/orm/5.2/userguide/html_single/Hibernate_User_Guide.html#mapping-model-pojo-equalshashcode"
$ws.Rows.Item(40).RowHeight = 75

$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D39").Value = "I read the following but did not come to a conclusion about making a change in unit tests (did not change):"

$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "Mon9am"

# Cursor ends up where the author last left it
$ws.Range("D58").Select() | Out-Null
